# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# (target/translation) files have now been produced for both content
# items, on both the zh-cn and de-de language sheets, and marks the
# overall status as handed back / in sync with en-US.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdName1 = "59a8820f-d6eb-4bca-82fe-1677208924ea.md"
$mdName2 = "f3fee816-963c-4e0f-ad79-c0a537ddebb3.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/901e0faa274ee8b0c0de8b83e4c576f8eec68186/e2e/59a8820f-d6eb-4bca-82fe-1677208924ea.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/901e0faa274ee8b0c0de8b83e4c576f8eec68186/e2e/f3fee816-963c-4e0f-ad79-c0a537ddebb3.md"

$handbackDateZhCn = "2016-08-16 04:58:32"
$handbackDateDeDe = "2016-08-16 04:58:39"

# ---------------------------------------------------------------------
# Overview sheet: refresh the cached status text for both items/langs
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: record the handback (target) files + datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew

$wsZhCn.Range("I2").Value = $mdName1
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZhCn.Range("I2").Style = "HyperLink"

$wsZhCn.Range("J2").Value = "59a8820f-d6eb-4bca-82fe-1677208924ea.43a1291d59df2bb4bf49dbde49470953867c741e.zh-cn.xlf"
$wsZhCn.Range("K2").Value = $handbackDateZhCn

$wsZhCn.Range("I3").Value = $mdName2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsZhCn.Range("I3").Style = "HyperLink"

$wsZhCn.Range("J3").Value = "f3fee816-963c-4e0f-ad79-c0a537ddebb3.3964448c1afcadce65023566408efed58c7363d1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = $handbackDateZhCn

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: record the handback (target) files + datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew

$wsDeDe.Range("I2").Value = $mdName1
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDeDe.Range("I2").Style = "HyperLink"

$wsDeDe.Range("J2").Value = "59a8820f-d6eb-4bca-82fe-1677208924ea.43a1291d59df2bb4bf49dbde49470953867c741e.de-de.xlf"
$wsDeDe.Range("K2").Value = $handbackDateDeDe

$wsDeDe.Range("I3").Value = $mdName2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsDeDe.Range("I3").Style = "HyperLink"

$wsDeDe.Range("J3").Value = "f3fee816-963c-4e0f-ad79-c0a537ddebb3.3964448c1afcadce65023566408efed58c7363d1.de-de.xlf"
$wsDeDe.Range("K3").Value = $handbackDateDeDe

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
